$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values: B2 and D2 updated, C2 and E2 cleared
$ws.Range("B2").Value = 1.6383776575388378
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 1.6710835076655888
$ws.Range("E2").ClearContents()

# Update row 3 values for columns B-E
$ws.Range("B3").Value = 1.2059266251779492
$ws.Range("C3").Value = -1.5116290300329904
$ws.Range("D3").Value = 1.2307636500082086
$ws.Range("E3").Value = -2.1901451881043488

# Update the selection to match new range B1:E3
$ws.Range("B1:E3").Select()
